$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("D2").Value = -0.03175
$ws.Range("E2").Value = -0.112
$ws.Range("F2").Value = 0.03555
$ws.Range("I2").Value = 0.001664375258189216
$ws.Range("J2").Value = 0.001206495860983133
$ws.Range("K2").Value = 12384.4
$ws.Range("L2").Value = 0.1171528168647687
$ws.Range("M2").Value = 10396.1972
$ws.Range("N2").Value = 0.05918451030244101
$ws.Range("O2").Value = 0.839459093698524
$ws.Range("P2").Value = 9904.6982
$ws.Range("Q2").Value = 0.05638645568020477
$ws.Range("R2").Value = 0.7997721488324021
$ws.Range("S2").Value = 491.4989999999997
$ws.Range("T2").Value = 0.04727680617678161
$ws.Range("U2").Value = 2590571.8
$ws.Range("V2").Value = 14.74786601646159
$ws.Range("W2").Value = 0.03604596813419601
$ws.Range("X2").Value = 0.2298621187821548
$ws.Range("Y2").Value = -0.1938161506479588
$ws.Range("Z2").Value = -0.8483005715892091
$ws.Range("AA2").Value = -0.000002867638158019965
$ws.Range("AB2").Value = 0.03774806730356173
$ws.Range("AC2").Value = -0.03744716853040073
$ws.Range("AD2").Value = 2063966.4
$ws.Range("AE2").Value = 6764.081974469654
$ws.Range("AF2").Value = 2070730.481974469
$ws.Range("AG2").Value = -519841.3180255303
$ws.Range("AH2").Value = 0.9218045105168546
$ws.Range("AI2").Value = 0.8130699358362655
$ws.Range("AJ2").Value = 1.510359115578928
$ws.Range("AK2").Value = 11.87749426912568
$ws.Range("AN2").Value = 1350.091839137602
$ws.Range("AP2").Value = -340.0411562478939

# Row 3
$ws.Range("B3").Value = "Sumitomo Mitsui Trust Holdings, Inc. (TSE:8309)"
$ws.Range("D3").Value = 0.0157
$ws.Range("E3").Value = -0.0328
$ws.Range("F3").Value = 0.0316
$ws.Range("I3").Value = -0.0002624567422960709
$ws.Range("J3").Value = -0.0001793935396297497
$ws.Range("K3").Value = 1299.7
$ws.Range("L3").Value = 0.124497107168857
$ws.Range("M3").Value = 531.9889999999999
$ws.Range("N3").Value = 0.04614275058113311
$ws.Range("O3").Value = 0.4093167654074016
$ws.Range("P3").Value = 531.79
$ws.Range("Q3").Value = 0.04612549006002151
$ws.Range("R3").Value = 0.4091636531507271
$ws.Range("S3").Value = 0.1989999999999554
$ws.Range("T3").Value = 0.0003740678848621972
$ws.Range("U3").Value = 165902.8
$ws.Range("V3").Value = 14.38979287374666
$ws.Range("W3").Value = 0.05255963862973703
$ws.Range("X3").Value = 0.246452745171988
$ws.Range("Y3").Value = -0.193893106542251
$ws.Range("Z3").Value = 0.8510659308541488
$ws.Range("AA3").Value = -0.0001526757297942135
$ws.Range("AB3").Value = 0.03556559105179079
$ws.Range("AC3").Value = -0.035718266781585
$ws.Range("AD3").Value = 135409.8
$ws.Range("AE3").Value = 191.6997170343703
$ws.Range("AF3").Value = 135601.4997170344
$ws.Range("AG3").Value = -30301.30028296562
$ws.Range("AH3").Value = 0.9216397392102853
$ws.Range("AI3").Value = 0.8442831577082168
$ws.Range("AJ3").Value = 1.614166759510759
$ws.Range("AK3").Value = 5.726518248962061
$ws.Range("AN3").Value = 3803.64606741573
$ws.Range("AP3").Value = -851.1601203080229

# Row 4
$ws.Range("B4").Value = "JAPAN POST BANK Co.,Ltd. (TSE:7182)"
$ws.Range("D4").Value = -0.0303
$ws.Range("E4").Value = -0.0679
$ws.Range("F4").Value = -0.0104
$ws.Range("I4").Value = 0.0001456823838313365
$ws.Range("J4").Value = 0.0001051520194834044
$ws.Range("K4").Value = 2396.7
$ws.Range("L4").Value = 0.1783059926347506
$ws.Range("M4").Value = 1779.7312
$ws.Range("N4").Value = 0.05785015846186351
$ws.Range("O4").Value = 0.7425757082655318
$ws.Range("P4").Value = 1776.9312
$ws.Range("Q4").Value = 0.05775914446846202
$ws.Range("R4").Value = 0.7414074352234322
$ws.Range("S4").Value = 2.799999999999955
$ws.Range("T4").Value = 0.001573271289507064
$ws.Range("U4").Value = 533870.1
$ws.Range("V4").Value = 17.35344634237514
$ws.Range("W4").Value = 0.02223360170989171
$ws.Range("X4").Value = 0.1543498118570249
$ws.Range("Y4").Value = -0.1321162101471332
$ws.Range("Z4").Value = -0.05454271201082445
$ws.Range("AA4").Value = -0.000005735276316039929
$ws.Range("AB4").Value = 0.03668659747653964
$ws.Range("AC4").Value = -0.03669233275285567
$ws.Range("AD4").Value = 204084.5
$ws.Range("AE4").Value = 6.009051188655453
$ws.Range("AF4").Value = 204090.5090511887
$ws.Range("AG4").Value = -329779.5909488113
$ws.Range("AH4").Value = 0.8690064132577449
$ws.Range("AI4").Value = 0.6603690961163682
$ws.Range("AJ4").Value = 1.102886111541663
$ws.Range("AK4").Value = 1.466893894603715
$ws.Range("AN4").Value = 64583.70253164557
$ws.Range("AP4").Value = -104360.6300470922

# Row 5
$ws.Range("B5").Value = "Mitsubishi UFJ Financial Group, Inc. (TSE:8306)"
$ws.Range("D5").Value = -0.0491
$ws.Range("E5").Value = -0.214
$ws.Range("F5").Value = 0.19
$ws.Range("I5").Value = 0.003818822811810438
$ws.Range("J5").Value = 0.002481492712225962
$ws.Range("K5").Value = 4421.2
$ws.Range("L5").Value = 0.1162919767163535
$ws.Range("M5").Value = 3527.7
$ws.Range("N5").Value = 0.06214842545694781
$ws.Range("O5").Value = 0.7979055460056094
$ws.Range("P5").Value = 3053.5
$ws.Range("Q5").Value = 0.05379431843206342
$ws.Range("R5").Value = 0.6906495973943726
$ws.Range("S5").Value = 474.1999999999998
$ws.Range("T5").Value = 0.1344218612693823
$ws.Range("U5").Value = 890455.1
$ws.Range("V5").Value = 15.68738339572781
$ws.Range("W5").Value = 0.02844100671848987
$ws.Range("X5").Value = 0.3204119800415812
$ws.Range("Y5").Value = -0.2919709733230913
$ws.Range("Z5").Value = 0.2595333470790933
$ws.Range("AA5").Value = 0.0006440301093563812
$ws.Range("AB5").Value = 0.03781435292744945
$ws.Range("AC5").Value = -0.03717032281809306
$ws.Range("AD5").Value = 897077.3
$ws.Range("AE5").Value = 4258.578062291548
$ws.Range("AF5").Value = 901335.8780622915
$ws.Range("AG5").Value = 10880.77806229156
$ws.Range("AH5").Value = 0.940755040088055
$ws.Range("AI5").Value = 0.8572955119370204
$ws.Range("AJ5").Value = 0.1608552745222021
$ws.Range("AK5").Value = 0.06761767644070908
$ws.Range("AN5").Value = 899.8668873507875
$ws.Range("AP5").Value = 10.91461336371909

# Row 6
$ws.Range("B6").Value = "Sumitomo Mitsui Financial Group, Inc. (TSE:8316)"
$ws.Range("D6").Value = -0.136
$ws.Range("E6").Value = -0.402
$ws.Range("F6").Value = 0.00866
$ws.Range("I6").Value = 0
$ws.Range("J6").Value = 0
$ws.Range("K6").Value = 479.7
$ws.Range("L6").Value = 0.02376116978066612
$ws.Range("M6").Value = 2583.8
$ws.Range("N6").Value = 0.06103858674339604
$ws.Range("O6").Value = 5.386283093600167
$ws.Range("P6").Value = 2583.8
$ws.Range("Q6").Value = 0.06103858674339604
$ws.Range("R6").Value = 5.386283093600167
$ws.Range("S6").Value = 0
$ws.Range("T6").Value = 0
$ws.Range("U6").Value = 600541.7
$ws.Range("V6").Value = 14.18694041662532
$ws.Range("W6").Value = 0.004763082092419877
$ws.Range("X6").Value = 0.2132714923923216
$ws.Range("Y6").Value = -0.2085084102999017
$ws.Range("Z6").Value = -0.4652959438373935
$ws.Range("AA6").Value = -0.0
$ws.Range("AB6").Value = 0.03772401424270838
$ws.Range("AC6").Value = -0.03772401424270838
$ws.Range("AD6").Value = 419677.6
$ws.Range("AE6").Value = 0
$ws.Range("AF6").Value = 419677.6
$ws.Range("AG6").Value = -180864.1
$ws.Range("AH6").Value = 0.9083769508852874
$ws.Range("AI6").Value = 0.7966365224475757
$ws.Range("AJ6").Value = 1.305562192538267
$ws.Range("AK6").Value = 2.453066467018764

# Row 7
$ws.Range("B7").Value = "Aozora Bank, Ltd. (TSE:8304)"
$ws.Range("D7").Value = 0.0151
$ws.Range("E7").Value = -0.123
$ws.Range("F7").Value = 0.0395
$ws.Range("I7").Value = 0.01110687203747664
$ws.Range("J7").Value = 0.006860542460807772
$ws.Range("K7").Value = 218.2
$ws.Range("L7").Value = 0.2045944678856071
$ws.Range("M7").Value = 152.877
$ws.Range("N7").Value = 0.07104610093874895
$ws.Range("O7").Value = 0.700627864344638
$ws.Range("P7").Value = 152.877
$ws.Range("Q7").Value = 0.07104610093874895
$ws.Range("R7").Value = 0.700627864344638
$ws.Range("S7").Value = 0
$ws.Range("T7").Value = 0
$ws.Range("U7").Value = 8094.8
$ws.Range("V7").Value = 3.761873780091086
$ws.Range("W7").Value = 0.0494034007290511
$ws.Range("X7").Value = 0.121434051612196
$ws.Range("Y7").Value = -0.07203065088314492
$ws.Range("Z7").Value = 0.1022475026445874
$ws.Range("AA7").Value = 0.0007014733334047469
$ws.Range("AB7").Value = 0.03845569972674823
$ws.Range("AC7").Value = -0.03775422639334348
$ws.Range("AD7").Value = 10298
$ws.Range("AE7").Value = 33.77260486015581
$ws.Range("AF7").Value = 10331.77260486016
$ws.Range("AG7").Value = 2236.972604860156
$ws.Range("AH7").Value = 0.8276294720982155
$ws.Range("AI7").Value = 0.6966818402263383
$ws.Range("AJ7").Value = 0.5097034652428603
$ws.Range("AK7").Value = 0.3321329290426704
$ws.Range("AN7").Value = 553.6559139784946
$ws.Range("AP7").Value = 120.2673443473202

# Row 8
$ws.Range("B8").Value = "Mizuho Financial Group, Inc. (TSE:8411)"
$ws.Range("D8").Value = -0.0332
$ws.Range("E8").Value = -0.101
$ws.Range("F8").Value = 0.0717
$ws.Range("I8").Value = 0.000873127762109782
$ws.Range("J8").Value = 0.0005907307630280142
$ws.Range("K8").Value = 3568.9
$ws.Range("L8").Value = 0.15821415588676
$ws.Range("M8").Value = 1820.1
$ws.Range("N8").Value = 0.05666774599300099
$ws.Range("O8").Value = 0.5099890722631623
$ws.Range("P8").Value = 1805.8
$ws.Range("Q8").Value = 0.05622252388009515
$ws.Range("R8").Value = 0.5059822354226793
$ws.Range("S8").Value = 14.29999999999995
$ws.Range("T8").Value = 0.007856711169715926
$ws.Range("U8").Value = 391707.3
$ws.Range("V8").Value = 12.19557704521962
$ws.Range("W8").Value = 0.04365092954990216
$ws.Range("X8").Value = 0.258714496099499
$ws.Range("Y8").Value = -0.2150635665495969
$ws.Range("Z8").Value = -5.680709168713153
$ws.Range("AA8").Value = -0.003355769661774157
$ws.Range("AB8").Value = 0.03777212036441507
$ws.Range("AC8").Value = -0.04112789002618923
$ws.Range("AD8").Value = 397419.2
$ws.Range("AE8").Value = 2274.022539094924
$ws.Range("AF8").Value = 399693.2225390949
$ws.Range("AG8").Value = 7985.922539094929
$ws.Range("AH8").Value = 0.925618560105996
$ws.Range("AI8").Value = 0.8255990113523146
$ws.Range("AJ8").Value = 0.1991267370397609
$ws.Range("AK8").Value = 0.08641106574131516
$ws.Range("AN8").Value = 837.5536354056902
$ws.Range("AP8").Value = 16.8301844870283

# Remove AN6 and AP6 cells entirely (no longer present after edit)
$ws.Range("AN6").ClearContents()
$ws.Range("AP6").ClearContents()
